$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OpenAccountTest")

# Update existing currency value from EURO to Dollar
$ws.Range("B2").Value = "Dollar"

# Add new column C: header "alerttext" and value "Account created successfully with account Number"
$ws.Range("C1").Value = "alerttext"
$ws.Range("C2").Value = "Account created successfully with account Number"

$ws.Range("C2").Select()
